$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 316
$ws.Range("H2").Value = 137.0601265822785
$ws.Range("I2").Value = 120.1679757883243
$ws.Range("O2").Value = 2.531645569620253
$ws.Range("P2").Value = 3.311683763132613
$ws.Range("V2").Value = 19.52215189873418
$ws.Range("W2").Value = 12.41847637009594
$ws.Range("Y2").Value = 18
$ws.Range("Z2").Value = 28
$ws.Range("AA2").Value = 40
$ws.Range("AC2").Value = 29.5
$ws.Range("AD2").Value = 17.23608249233391
$ws.Range("AJ2").Value = 1.882911392405063
$ws.Range("AK2").Value = 1.152867675771185
$ws.Range("AM2").Value = 0.5
$ws.Range("AQ2").Value = 0.6023722280384506
$ws.Range("AR2").Value = 0.3373627583929834
$ws.Range("AU2").Value = 38.90909090909091
$ws.Range("AX2").Value = 43.01886680478453
$ws.Range("AY2").Value = 14.97922458263611
$ws.Range("B3").Value = 316
$ws.Range("H3").Value = 137.0601265822785
$ws.Range("I3").Value = 120.1679757883243
$ws.Range("R3").Value = 18
$ws.Range("S3").Value = 29.5
$ws.Range("T3").Value = 57.25
$ws.Range("V3").Value = 43.62341772151899
$ws.Range("W3").Value = 38.24732324730301
$ws.Range("Y3").Value = 46
$ws.Range("Z3").Value = 65
$ws.Range("AC3").Value = 90.98417721518987
$ws.Range("AD3").Value = 83.66799317647113
$ws.Range("AJ3").Value = 4.218354430379747
$ws.Range("AK3").Value = 3.993619333229222
$ws.Range("AX3").Value = 137.0601265822785
$ws.Range("AY3").Value = 120.1679757883243
$ws.Range("B4").Value = 316
$ws.Range("H4").Value = 137.0601265822785
$ws.Range("I4").Value = 120.1679757883243
$ws.Range("R4").Value = 26
$ws.Range("S4").Value = 41
$ws.Range("T4").Value = 57.25
$ws.Range("V4").Value = 55.30696202531646
$ws.Range("W4").Value = 55.42210964128501
$ws.Range("Y4").Value = 43
$ws.Range("Z4").Value = 57
$ws.Range("AC4").Value = 86.1613924050633
$ws.Range("AD4").Value = 88.13580374369498
$ws.Range("AJ4").Value = 3.607594936708861
$ws.Range("AK4").Value = 3.430647017577352
$ws.Range("AX4").Value = 137.0601265822785
$ws.Range("AY4").Value = 120.1679757883243
$ws.Range("B5").Value = 316
$ws.Range("H5").Value = 137.0601265822785
$ws.Range("I5").Value = 120.1679757883243
$ws.Range("O5").Value = 3.522151898734177
$ws.Range("P5").Value = 4.594974804051639
$ws.Range("R5").Value = 9.75
$ws.Range("U5").Value = 58
$ws.Range("V5").Value = 17.66455696202532
$ws.Range("W5").Value = 11.07097492702609
$ws.Range("Y5").Value = 19
$ws.Range("Z5").Value = 30
$ws.Range("AA5").Value = 40
$ws.Range("AC5").Value = 30.17721518987342
$ws.Range("AD5").Value = 15.16471216599101
$ws.Range("AJ5").Value = 1.620253164556962
$ws.Range("AK5").Value = 0.9092651485168747
$ws.Range("AQ5").Value = 0.8354509691951397
$ws.Range("AR5").Value = 0.2479190750366284
$ws.Range("AW5").Value = 74
$ws.Range("AX5").Value = 33.1570297554705
$ws.Range("AY5").Value = 6.900392719654174
